$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the fill style previously applied to A9:C10 (back to default/no style)
$ws.Range("A9:C10").Style = "Normal"

# Add new rows 30-32 with additional country data
$ws.Range("A30").Value = 2020
$ws.Range("B30").Value = "NIC"
$ws.Range("C30").Value = 9.3699999999999992

$ws.Range("A31").Value = 2022
$ws.Range("B31").Value = "NIC"
$ws.Range("C31").Value = 11.4

$ws.Range("A32").Value = 2015
$ws.Range("B32").Value = "DOM"
$ws.Range("C32").Value = 16.48

# Update the view: scroll and select to match the end-state
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("E30").Select()
